# Insert a new "Match ID" column at the very left of the data table.
# This shifts every existing column (Player ID, Player, #, Nation, Pos, Age,
# 90s, Att, Live, Dead, FK, TB, Sw, Crs, TI, CK, In, Out, Str, Cmp, Off,
# Blocks) one column to the right (A->B, B->C, ... V->W) while Excel keeps
# styles/merged-cell ranges/formula refs consistent automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the right, inserting new blank column A.
$ws.Columns("A").Insert()

# Header label for the new column (row 2 holds the human-readable headers).
$ws.Range("A2").Value2 = "Match ID"

# Match the bold "header" styling used by the other id-style columns
# (fontId=1 / no border / no fill) for the header cell and the data column
# (rows 2 through 19, which also covers the blank hidden spacer row 3).
$ws.Range("A2:A19").Font.Bold = $true

# Fill in the constant Match ID value for every visible + hidden data row.
$ws.Range("A4:A19").Value2 = 23
$ws.Range("A20").Value2 = 23

# Re-fit row 20's height so writing into it doesn't leave a stray explicit
# row-height override (it's a hidden helper/summary row).
$ws.Rows(20).AutoFit()

# Restore the sheet's remembered selection to the new data column.
$ws.Range("A2:A19").Select()
